$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 4 and row 5 for columns D, M, N, O, P, S
$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cell4 = $ws.Range($col + "4")
    $cell5 = $ws.Range($col + "5")
    $v4 = $cell4.Value2
    $v5 = $cell5.Value2
    $cell4.Value2 = $v5
    $cell5.Value2 = $v4
}
